# Fruta / hortaliza, semanal
# Inserts one new week's worth of price data (3 rows: Especial / Primera /
# Segunda) for Femacal de La Calera - Piña, right before the current row 471,
# pushing the existing data down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 471 (existing rows 471:569 shift down to 474:572).
$ws.Rows("471:473").Insert()

# Columns that are constant across every data row in this sheet.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100108
$producto   = "Tropicales y subtropicales"
$categoriaId = 100108005
$categoria   = "Piña"
$variedad    = "Caramelo"
$origen      = "Ecuador"
$fecha       = 44637

# Row 471 - Especial
$ws.Cells.Item(471, 1).Value = $mercadoId
$ws.Cells.Item(471, 2).Value = $mercado
$ws.Cells.Item(471, 3).Value = $region
$ws.Cells.Item(471, 4).Value = $fecha
$ws.Cells.Item(471, 5).Value = $codreg
$ws.Cells.Item(471, 6).Value = $tipo
$ws.Cells.Item(471, 7).Value = $productoId
$ws.Cells.Item(471, 8).Value = $producto
$ws.Cells.Item(471, 9).Value = $categoriaId
$ws.Cells.Item(471, 10).Value = $categoria
$ws.Cells.Item(471, 11).Value = $variedad
$ws.Cells.Item(471, 12).Value = "Especial"
$ws.Cells.Item(471, 13).Value = 108
$ws.Cells.Item(471, 14).Value = 18000
$ws.Cells.Item(471, 15).Value = 18000
$ws.Cells.Item(471, 16).Value = 18000
$ws.Cells.Item(471, 17).Value = "$/caja 10 unidades"
$ws.Cells.Item(471, 18).Value = $origen
$ws.Cells.Item(471, 19).Value = 1800
$ws.Cells.Item(471, 20).Value = 10

# Row 472 - Primera
$ws.Cells.Item(472, 1).Value = $mercadoId
$ws.Cells.Item(472, 2).Value = $mercado
$ws.Cells.Item(472, 3).Value = $region
$ws.Cells.Item(472, 4).Value = $fecha
$ws.Cells.Item(472, 5).Value = $codreg
$ws.Cells.Item(472, 6).Value = $tipo
$ws.Cells.Item(472, 7).Value = $productoId
$ws.Cells.Item(472, 8).Value = $producto
$ws.Cells.Item(472, 9).Value = $categoriaId
$ws.Cells.Item(472, 10).Value = $categoria
$ws.Cells.Item(472, 11).Value = $variedad
$ws.Cells.Item(472, 12).Value = "Primera"
$ws.Cells.Item(472, 13).Value = 162
$ws.Cells.Item(472, 14).Value = 18000
$ws.Cells.Item(472, 15).Value = 18000
$ws.Cells.Item(472, 16).Value = 18000
$ws.Cells.Item(472, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(472, 18).Value = $origen
$ws.Cells.Item(472, 19).Value = 1500
$ws.Cells.Item(472, 20).Value = 12

# Row 473 - Segunda
$ws.Cells.Item(473, 1).Value = $mercadoId
$ws.Cells.Item(473, 2).Value = $mercado
$ws.Cells.Item(473, 3).Value = $region
$ws.Cells.Item(473, 4).Value = $fecha
$ws.Cells.Item(473, 5).Value = $codreg
$ws.Cells.Item(473, 6).Value = $tipo
$ws.Cells.Item(473, 7).Value = $productoId
$ws.Cells.Item(473, 8).Value = $producto
$ws.Cells.Item(473, 9).Value = $categoriaId
$ws.Cells.Item(473, 10).Value = $categoria
$ws.Cells.Item(473, 11).Value = $variedad
$ws.Cells.Item(473, 12).Value = "Segunda"
$ws.Cells.Item(473, 13).Value = 162
$ws.Cells.Item(473, 14).Value = 18000
$ws.Cells.Item(473, 15).Value = 18000
$ws.Cells.Item(473, 16).Value = 18000
$ws.Cells.Item(473, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(473, 18).Value = $origen
$ws.Cells.Item(473, 19).Value = 1286
$ws.Cells.Item(473, 20).Value = 14
